$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "275.08"
Set-TextValue "G2" "16"
Set-TextValue "D3" "23.02"
Set-TextValue "G3" "16"
Set-TextValue "D4" "6.308"
Set-TextValue "G4" "16"
Set-TextValue "D5" "0.06235"
Set-TextValue "G5" "16"
Set-TextValue "D6" "3.647"
Set-TextValue "G6" "16"
Set-TextValue "D7" "6.654"
Set-TextValue "G7" "16"
Set-TextValue "D8" "1.394"
Set-TextValue "G8" "16"
Set-TextValue "D9" "0.8325"
Set-TextValue "G9" "16"
Set-TextValue "D10" "0.01382"
Set-TextValue "G10" "16"
Set-TextValue "D11" "0.1594"
Set-TextValue "G11" "16"
Set-TextValue "D12" "0.08391"
Set-TextValue "G12" "16"
Set-TextValue "D13" "0.03535"
Set-TextValue "G13" "16"
Set-TextValue "D14" "0.03185"
Set-TextValue "G14" "16"
Set-TextValue "B15" "MCDex"
Set-TextValue "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "4.061"
Set-TextValue "E15" "14MCDexMCB"
Set-TextValue "G15" "16"
Set-TextValue "B16" "BitMartToken"
Set-TextValue "C16" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D16" "0.09300"
Set-TextValue "E16" "15BitMartTokenBMX"
Set-TextValue "G16" "16"
Set-TextValue "B17" "BitForexToken"
Set-TextValue "C17" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001635"
Set-TextValue "E17" "16BitForexTokenBF"
Set-TextValue "G17" "16"
Set-TextValue "B18" "CoinExToken"
Set-TextValue "C18" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04729"
Set-TextValue "E18" "17CoinExTokenCET"
Set-TextValue "G18" "16"
Set-TextValue "B19" "TigerCash"
Set-TextValue "C19" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D19" "0.006272"
Set-TextValue "E19" "18TigerCashTCH"
Set-TextValue "G19" "16"
Set-TextValue "B20" "HotbitToken"
Set-TextValue "C20" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D20" "0.005709"
Set-TextValue "E20" "19HotbitTokenHTB"
Set-TextValue "G20" "16"
Set-TextValue "B21" "BitKan"
Set-TextValue "C21" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D21" "0.001079"
Set-TextValue "E21" "20BitKanKAN"
Set-TextValue "G21" "16"
Set-TextValue "B22" "NitroEx"
Set-TextValue "C22" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D22" "0.0001500"
Set-TextValue "E22" "21NitroExNTX"
Set-TextValue "G22" "16"
Set-TextValue "B23" "LEO"
Set-TextValue "C23" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D23" "3.733"
Set-TextValue "E23" "22LEOLEO"
Set-TextValue "G23" "16"
Set-TextValue "B24" "BTSEToken"
Set-TextValue "C24" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D24" "2.414"
Set-TextValue "E24" "23BTSETokenBTSE"
Set-TextValue "G24" "16"
Set-TextValue "B25" "BitpandaEcosystemToken"
Set-TextValue "C25" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D25" "0.3358"
Set-TextValue "E25" "24BitpandaEcosystemTokenBEST"
Set-TextValue "G25" "16"
Set-TextValue "B26" "ProBitToken"
Set-TextValue "C26" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D26" "0.1261"
Set-TextValue "E26" "25ProBitTokenPROB"
Set-TextValue "G26" "16"
Set-TextValue "B27" "UpBots"
Set-TextValue "C27" "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue "D27" "0.0002703"
Set-TextValue "E27" "26UpBotsUBXT"
Set-TextValue "G27" "16"
Set-TextValue "B28" "Spectre.aiUtilityToken"
Set-TextValue "C28" "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
Set-TextValue "D28" "--"
Set-TextValue "E28" "27Spectre.aiUtilityTokenSXUT"
Set-TextValue "G28" "16"
Set-TextValue "B29" "LegolasExchange"
Set-TextValue "C29" "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
Set-TextValue "E29" "28LegolasExchangeLGO"
Set-TextValue "G29" "16"
Set-TextValue "B30" "BitZToken"
Set-TextValue "C30" "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
Set-TextValue "E30" "29BitZTokenBZ"
Set-TextValue "G30" "16"
Set-TextValue "B31" "Birake"
Set-TextValue "C31" "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
Set-TextValue "E31" "30BirakeBIR"
Set-TextValue "G31" "16"
Set-TextValue "B32" "ZBToken"
Set-TextValue "C32" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "E32" "31ZBTokenZB"
Set-TextValue "G32" "16"
Set-TextValue "B33" "NashExchange"
Set-TextValue "C33" "https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"
Set-TextValue "E33" "32NashExchangeNEX"
Set-TextValue "G33" "16"
Set-TextValue "B34" "AAXToken"
Set-TextValue "C34" "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
Set-TextValue "E34" "33AAXTokenAAB"
Set-TextValue "G34" "16"
Set-TextValue "G35" "16"
Set-TextValue "G36" "16"
Set-TextValue "G37" "16"
Set-TextValue "G38" "16"
Set-TextValue "G39" "16"
Set-TextValue "D40" "0.04735"
Set-TextValue "G40" "16"
Set-TextValue "D41" "0.006997"
Set-TextValue "G41" "16"
Set-TextValue "D42" "0.003899"
Set-TextValue "E42" "41CEJICEJI"
Set-TextValue "G42" "16"
Set-TextValue "D43" "0.1168"
Set-TextValue "G43" "16"
Set-TextValue "D44" "0.01181"
Set-TextValue "G44" "16"
Set-TextValue "D45" "0.00006239"
Set-TextValue "G45" "16"
Set-TextValue "D46" "0.0009899"
Set-TextValue "G46" "16"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "G47" "16"
Set-TextValue "D48" "0.7820"
Set-TextValue "E48" "47CoinbaseStockTokenCOINWorstin24h"
Set-TextValue "G48" "16"
Set-TextValue "D49" "0.002215"
Set-TextValue "G49" "16"
Set-TextValue "D50" "0.00002399"
Set-TextValue "G50" "16"
Set-TextValue "D51" "0.01240"
Set-TextValue "G51" "16"
